$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.920.41"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").Value = "3.131.78"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.23%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "3.121.46"
$ws.Range("E8").Value = "  +0.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("E10").Value = "  -0.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.35"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.467"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000252"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("D15").Value = "3.650.67"
$ws.Range("E15").Value = "  +0.66%  "
$ws.Range("E16").Value = "  +2.83%  "
$ws.Range("D17").Value = "63.992.79"
$ws.Range("E17").Value = "  +0.89%  "
$ws.Range("D18").Value = "3.134.32"
$ws.Range("E18").Value = "  +0.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "479.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.708"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.62"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  -1.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "27.04"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.83%  "
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.110"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.64"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.60%  "
$ws.Range("E35").Value = "  -1.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.92%  "
$ws.Range("D37").Value = "0.0₃0751"
$ws.Range("E37").Value = "  -3.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.61"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "437.37"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0393"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("E42").Value = "  +0.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.25"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").Value = "2.871.84"
$ws.Range("E44").Value = "  +0.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.259"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.84%  "
$ws.Range("E46").Value = "  +2.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.80"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.92%  "
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "121.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.53%  "
